$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 3845.8276
$ws.Range("I9").Value = 5698.6113
$ws.Range("K9").Value = 5698.6113
$ws.Range("M9").Value = -5529.6113

# Row 28
$ws.Range("H28").Value = 2437.5386
$ws.Range("I28").Value = 1974
$ws.Range("J28").Value = 8000
$ws.Range("K28").Value = 1974
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = -1489
$ws.Range("N28").Value = -8970

# Row 40
$ws.Range("H40").Value = 3839155
$ws.Range("I40").Value = 4348.7915
$ws.Range("K40").Value = 4348.7915
$ws.Range("M40").Value = -4173.7915

# Row 48
$ws.Range("H48").Value = 9194.5
$ws.Range("I48").Value = 8926.666999999999
$ws.Range("K48").Value = 26780.001
$ws.Range("M48").Value = -26488.001

# Row 56
$ws.Range("H56").Value = 9194.5
$ws.Range("I56").Value = 8926.666999999999
$ws.Range("K56").Value = 26780.001
$ws.Range("M56").Value = -26246.001

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

# Row 97
$ws.Range("H97").Value = 3637
$ws.Range("J97").Value = 3637
$ws.Range("L97").Value = 10911
$ws.Range("N97").Value = -11903

# Row 112
$ws.Range("H112").Value = 69338.07000000001
$ws.Range("I112").Value = 85340.086
$ws.Range("J112").Value = 58670.055
$ws.Range("K112").Value = 256020.258
$ws.Range("L112").Value = 176010.165
$ws.Range("M112").Value = -254912.258
$ws.Range("N112").Value = -178226.165

# Row 141
$ws.Range("H141").Value = 3819
$ws.Range("I141").Value = 3819
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11457
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -6277
$ws.Range("N141").ClearContents()


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 97
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504

# Row 102
$ws.Range("H102").Value = 30303596
$ws.Range("I102").Value = 30303596
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 30303596
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -30301974
$ws.Range("N102").ClearContents()

# Row 110
$ws.Range("H110").Value = 34812.633
$ws.Range("I110").Value = 41113
$ws.Range("J110").Value = 3310.8
$ws.Range("K110").Value = 41113
$ws.Range("L110").Value = 3310.8
$ws.Range("M110").Value = -39068
$ws.Range("N110").Value = -7400.8


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 1887.3334
$ws.Range("I99").Value = 1642
$ws.Range("K99").Value = 1642
$ws.Range("M99").Value = -144

# Row 107
$ws.Range("H107").Value = 60533.117
$ws.Range("J107").Value = 144981.28
$ws.Range("L107").Value = 144981.28
$ws.Range("N107").Value = -148821.28


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 17919.4
$ws.Range("I31").Value = 13912.571
$ws.Range("K31").Value = 13912.571
$ws.Range("M31").Value = -13617.571

# Row 34
$ws.Range("H34").Value = 17919.4
$ws.Range("I34").Value = 13912.571
$ws.Range("K34").Value = 13912.571
$ws.Range("M34").Value = -13710.571

# Row 141
$ws.Range("H141").Value = 455997.6
$ws.Range("J141").Value = 455997.6
$ws.Range("L141").Value = 455997.6
$ws.Range("N141").Value = -466357.6


# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 271.3125
$ws.Range("I2").Value = 257.14285
$ws.Range("J2").Value = 282.33334
$ws.Range("K2").Value = 1542.8571
$ws.Range("L2").Value = 1694.00004
$ws.Range("M2").Value = -1429.8571
$ws.Range("N2").Value = -1920.00004

# Row 5
$ws.Range("H5").Value = 125545.375
$ws.Range("J5").Value = 504.25
$ws.Range("L5").Value = 1512.75
$ws.Range("N5").Value = -1736.75

# Row 114
$ws.Range("H114").Value = 67835.336
$ws.Range("J114").Value = 2182.25
$ws.Range("L114").Value = 6546.75
$ws.Range("N114").Value = -13054.75

# Row 131
$ws.Range("H131").Value = 1703.7
$ws.Range("I131").Value = 947.7646999999999
$ws.Range("J131").Value = 2692.2307
$ws.Range("K131").Value = 2843.2941
$ws.Range("L131").Value = 8076.6921
$ws.Range("M131").Value = 2196.7059
$ws.Range("N131").Value = -18156.6921

# Row 135
$ws.Range("H135").Value = 125545.375
$ws.Range("J135").Value = 504.25
$ws.Range("L135").Value = 4538.25
$ws.Range("N135").Value = -9608.25


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 57
$ws.Range("H57").Value = 56799.8
$ws.Range("J57").Value = 63500
$ws.Range("L57").Value = 63500
$ws.Range("N57").Value = -65140

# Row 80
$ws.Range("H80").Value = 3248.5
$ws.Range("I80").Value = 3248.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3248.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2250.5
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 3248.5
$ws.Range("I83").Value = 3248.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16242.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11250.5
$ws.Range("N83").ClearContents()

# Row 102
$ws.Range("H102").Value = 4209.1113
$ws.Range("I102").Value = 1298.0625
$ws.Range("K102").Value = 1298.0625
$ws.Range("M102").Value = 323.9375

# Row 113
$ws.Range("H113").Value = 74984.14
$ws.Range("I113").Value = 113534.89
$ws.Range("K113").Value = 113534.89
$ws.Range("M113").Value = -111364.89


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 82
$ws.Range("H82").Value = 997.8
$ws.Range("I82").Value = 997.5
$ws.Range("K82").Value = 997.5
$ws.Range("M82").Value = -636.5

# Row 85
$ws.Range("H85").Value = 997.8
$ws.Range("I85").Value = 997.5
$ws.Range("K85").Value = 997.5
$ws.Range("M85").Value = 250.5

# Row 122
$ws.Range("H122").Value = 16666
$ws.Range("I122").Value = 16666
$ws.Range("K122").Value = 49998
$ws.Range("M122").Value = -47548

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()


# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 4
$ws.Range("H4").Value = 150002500
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 300000000
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 300000000
$ws.Range("M4").Value = -4887
$ws.Range("N4").Value = -300000226

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 126
$ws.Range("H126").Value = 1198.1428
$ws.Range("I126").Value = 1064.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3193.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -723.5
$ws.Range("N126").Value = -10940

